$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{
        B = 0.00158971523361902
        C = 0.0000691180536356096
        D = 0
        E = 0
        F = 0
        G = 0.999792645839093
        H = 0.442977605750622
        I = 0.999308819463644
        J = 0.000207354160906829
        K = 0
        L = 0.996474979264584
        M = 0.0000691180536356096
        N = 0.000552944429084877
        O = 0.000276472214542438
        P = 0.000207354160906829
        Q = 0.000483826375449267
        R = 0.000207354160906829
        S = 0.000898534697262925
        T = 0.000345590268178048
        U = 0.000276472214542438
        V = 0.996405861210948
        W = 0.915883328725463
        X = 0.000967652750898535
    }
    3 = @{
        B = 0.997649986176389
        C = 0.999723527785458
        D = 0.999239701410008
        E = 0.000276472214542438
        F = 0.999792645839093
        G = 0
        H = 0.0000691180536356096
        I = 0.000276472214542438
        J = 0.000483826375449267
        K = 0.99149847940282
        L = 0.00138236107271219
        M = 0.998755875034559
        N = 0.000276472214542438
        O = 0.998963229195466
        P = 0.000760298589991706
        Q = 0.000829416643627315
        R = 0.000552944429084877
        S = 0
        T = 0.000138236107271219
        U = 0.999032347249102
        V = 0.000483826375449267
        W = 0.000138236107271219
        X = 0.998686756980923
    }
    4 = @{
        B = 0.0000691180536356096
        C = 0
        D = 0
        E = 0
        F = 0.0000691180536356096
        G = 0.0000691180536356096
        H = 0.548106165330384
        I = 0.000414708321813658
        J = 0.000345590268178048
        K = 0
        L = 0.000967652750898535
        M = 0.000552944429084877
        N = 0
        O = 0.0000691180536356096
        P = 0
        Q = 0.0000691180536356096
        R = 0.000483826375449267
        S = 0.999032347249102
        T = 0.999447055570915
        U = 0.000345590268178048
        V = 0.00304119435996682
        W = 0.080729886646392
        X = 0.000138236107271219
    }
    5 = @{
        B = 0.000414708321813658
        C = 0.000138236107271219
        D = 0.000552944429084877
        E = 0.999516173624551
        F = 0.000138236107271219
        G = 0.0000691180536356096
        H = 0.000138236107271219
        I = 0
        J = 0.998755875034559
        K = 0.00836328448990876
        L = 0.000829416643627315
        M = 0.000207354160906829
        N = 0.999101465302737
        O = 0.000207354160906829
        P = 0.998963229195466
        Q = 0.998548520873652
        R = 0.998548520873652
        S = 0.0000691180536356096
        T = 0.0000691180536356096
        U = 0.000207354160906829
        V = 0
        W = 0.0000691180536356096
        X = 0
    }
}

foreach ($rowNum in $data.Keys) {
    $rowData = $data[$rowNum]
    foreach ($col in $rowData.Keys) {
        $addr = "$col$rowNum"
        $ws.Range($addr).Value = $rowData[$col]
    }
}
